# Updated capital structure database for Hungary / Investments & Asset Management
# - AKKO Invest Nyrt. (row 6) removed from the dataset
# - refreshed metrics for the remaining four companies (rows 2-5)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: Hungary / company #2 - refreshed metrics ---
$ws.Range("B2").Value = "3"
$ws.Range("G2").Value = 0.02889833641404806
$ws.Range("H2").Value = 0.02889833641404806
$ws.Range("I2").Value = -0.007390018484288353
$ws.Range("J2").Value = -0.006158348736906961
$ws.Range("K2").Value = -3.316
$ws.Range("L2").Value = -0.01225878003696858
$ws.Range("M2").Value = 3.46
$ws.Range("N2").Value = 0.03002950876583926
$ws.Range("O2").Value = -1.043425814234017
$ws.Range("P2").Value = 1.04
$ws.Range("Q2").Value = 0.009026210727304288
$ws.Range("R2").Value = -0.3136308805790109
$ws.Range("S2").Value = 2.42
$ws.Range("T2").Value = 0.6994219653179191
$ws.Range("U2").Value = 28.239
$ws.Range("V2").Value = 0.2450876583926402
$ws.Range("W2").Value = -0.02985781990521327
$ws.Range("X2").Value = 0.07054030257711663
$ws.Range("Y2").Value = -0.1003981224823299
$ws.Range("Z2").Value = 2.131582887447695
$ws.Range("AA2").Value = 0.01757402101241643
$ws.Range("AB2").Value = 0.05857512243026935
$ws.Range("AC2").Value = -0.04100110141785292
$ws.Range("AD2").Value = 84.5
$ws.Range("AF2").Value = 84.5
$ws.Range("AG2").Value = 56.261
$ws.Range("AH2").Value = 0.4230923292609653
$ws.Range("AI2").Value = 0.4552041415496334
$ws.Range("AJ2").Value = 0.3280888261673305
$ws.Range("AK2").Value = 0.3574578123411609
$ws.Range("AL2").Value = 1.792
$ws.Range("AM2").Value = 1.634
$ws.Range("AN2").Value = 41.21951219512196
$ws.Range("AO2").Value = -1.115513392857143
$ws.Range("AP2").Value = 27.44439024390244
$ws.Range("AQ2").Value = -1.223378212974296

# --- Row 3: Elso Hazai Energia-portfolió - refreshed metrics ---
$ws.Range("K3").Value = -0.036
$ws.Range("U3").Value = 0.039
$ws.Range("V3").Value = 0.01547619047619048
$ws.Range("W3").Value = 1.5
$ws.Range("X3").Value = 0.04696437217823866
$ws.Range("Y3").Value = 1.453035627821761
$ws.Range("AA3").Value = 1.625
$ws.Range("AB3").Value = 0.04696437217823866
$ws.Range("AC3").Value = 1.578035627821761
$ws.Range("AG3").Value = -0.039
$ws.Range("AJ3").Value = -0.01571946795646917
$ws.Range("AK3").Value = 4.875
$ws.Range("AL3").Value = 0
$ws.Range("AM3").Value = 0
# cells no longer populated for this company
$ws.Range("AN3").ClearContents()
$ws.Range("AO3").ClearContents()
$ws.Range("AP3").ClearContents()
$ws.Range("AQ3").ClearContents()

# --- Row 4: AutoWallis - refreshed metrics ---
$ws.Range("G4").Value = 0.02998102466793169
$ws.Range("H4").Value = 0.02998102466793169
$ws.Range("I4").Value = 0.00698292220113852
$ws.Range("J4").Value = 0.00349146110056926
$ws.Range("K4").Value = -0.76
$ws.Range("L4").Value = -0.002884250474383302
$ws.Range("M4").Value = 2.42
$ws.Range("N4").Value = 0.02746878547105562
$ws.Range("O4").Value = -3.184210526315789
$ws.Range("S4").Value = 2.42
$ws.Range("U4").Value = 14.2
$ws.Range("V4").Value = 0.1611804767309875
$ws.Range("W4").Value = -0.04606060606060606
$ws.Range("X4").Value = 0.07054030257711663
$ws.Range("Y4").Value = -0.1166009086377227
$ws.Range("Z4").Value = 5.033428844317096
$ws.Range("AA4").Value = 0.01757402101241643
$ws.Range("AB4").Value = 0.05857512243026935
$ws.Range("AC4").Value = -0.04100110141785292
$ws.Range("AD4").Value = 60.6
$ws.Range("AF4").Value = 60.6
$ws.Range("AG4").Value = 46.40000000000001
$ws.Range("AH4").Value = 0.4075319435104237
$ws.Range("AI4").Value = 0.7444717444717445
$ws.Range("AJ4").Value = 0.3449814126394052
$ws.Range("AK4").Value = 0.6904761904761906
$ws.Range("AL4").Value = 0.582
$ws.Range("AM4").Value = 0.5509999999999999
$ws.Range("AN4").Value = 12.41803278688525
$ws.Range("AO4").Value = 3.161512027491409
$ws.Range("AP4").Value = 9.508196721311476
$ws.Range("AQ4").Value = 3.339382940108893
$ws.Range("T4").Value = 1

# --- Row 5: Forras nyRt. - refreshed metrics ---
$ws.Range("G5").Value = -0.006285714285714285
$ws.Range("H5").Value = -0.006285714285714285
$ws.Range("I5").Value = -0.5428571428571428
$ws.Range("J5").Value = -0.5428571428571428
$ws.Range("K5").Value = -2.52
$ws.Range("L5").Value = -0.36
$ws.Range("M5").Value = 1.04
$ws.Range("N5").Value = 0.04227642276422764
$ws.Range("O5").Value = -0.4126984126984127
$ws.Range("P5").Value = 1.04
$ws.Range("Q5").Value = 0.04227642276422764
$ws.Range("R5").Value = -0.4126984126984127
$ws.Range("U5").Value = 14
$ws.Range("V5").Value = 0.5691056910569106
$ws.Range("W5").Value = -0.02985781990521327
$ws.Range("X5").Value = 0.08026365797110502
$ws.Range("Y5").Value = -0.1101214778763183
$ws.Range("Z5").Value = 0.09386523633925578
$ws.Range("AA5").Value = -0.05095541401273885
$ws.Range("AB5").Value = 0.06100397397579552
$ws.Range("AC5").Value = -0.1119593879885344
$ws.Range("AD5").Value = 23.9
$ws.Range("AF5").Value = 23.9
$ws.Range("AG5").Value = 9.899999999999999
$ws.Range("AH5").Value = 0.4927835051546391
$ws.Range("AI5").Value = 0.2293666026871401
$ws.Range("AJ5").Value = 0.2869565217391304
$ws.Range("AK5").Value = 0.1097560975609756
$ws.Range("AL5").Value = 1.21
$ws.Range("AM5").Value = 1.083
$ws.Range("AN5").Value = -8.445229681978798
$ws.Range("AO5").Value = -3.140495867768595
$ws.Range("AP5").Value = -3.498233215547703
$ws.Range("AQ5").Value = -3.508771929824561

# --- Row 6: AKKO Invest Nyrt. removed from dataset entirely ---
$ws.Rows.Item(6).Delete()
